$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 47 (ALC)
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").ClearContents()

# Row 58 (ALC)
$ws.Range("H58").Value = 679.8
$ws.Range("J58").Value = 1003.4
$ws.Range("L58").Value = 3010.2
$ws.Range("N58").Value = -3310.2

# Row 62 (ALC)
$ws.Range("H62").Value = 18523184
$ws.Range("I62").Value = 55560556
$ws.Range("J62").Value = 4499
$ws.Range("K62").Value = 55560556
$ws.Range("L62").Value = 4499
$ws.Range("M62").Value = -55559932
$ws.Range("N62").Value = -5747

# Row 65 (ALC)
$ws.Range("H65").Value = 18523184
$ws.Range("I65").Value = 55560556
$ws.Range("J65").Value = 4499
$ws.Range("K65").Value = 277802780
$ws.Range("L65").Value = 22495
$ws.Range("M65").Value = -277799660
$ws.Range("N65").Value = -28735

# Row 74 (ALC)
$ws.Range("H74").Value = 3714.7144
$ws.Range("I74").Value = 3334.3333
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3334.3333
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2398.3333
$ws.Range("N74").Value = -5872

# Row 76 (ALC)
$ws.Range("H76").Value = 5292.7144
$ws.Range("I76").Value = 4866.6665
$ws.Range("J76").Value = 5612.25
$ws.Range("K76").Value = 4866.6665
$ws.Range("L76").Value = 5612.25
$ws.Range("M76").Value = -4551.6665
$ws.Range("N76").Value = -6242.25

# Row 77 (ALC)
$ws.Range("H77").Value = 3714.7144
$ws.Range("I77").Value = 3334.3333
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 16671.6665
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -11991.6665
$ws.Range("N77").Value = -29360

# Row 79 (ALC)
$ws.Range("H79").Value = 5292.7144
$ws.Range("I79").Value = 4866.6665
$ws.Range("J79").Value = 5612.25
$ws.Range("K79").Value = 4866.6665
$ws.Range("L79").Value = 5612.25
$ws.Range("M79").Value = -3774.6665
$ws.Range("N79").Value = -7796.25

# Row 133 (ALC)
$ws.Range("H133").Value = 35320
$ws.Range("J133").Value = 35320
$ws.Range("L133").Value = 35320
$ws.Range("N133").Value = -45440

# Row 137 (ALC)
$ws.Range("H137").Value = 1549.9375
$ws.Range("I137").Value = 1262.2069
$ws.Range("J137").Value = 1989.1052
$ws.Range("K137").Value = 3786.620699999999
$ws.Range("L137").Value = 5967.3156
$ws.Range("M137").Value = -1236.620699999999
$ws.Range("N137").Value = -11067.3156

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 755.8823
$ws.Range("J2").Value = 1830
$ws.Range("K2").Value = 755.8823
$ws.Range("L2").Value = 1830
$ws.Range("M2").Value = -642.8823
$ws.Range("N2").Value = -2056

# Row 32 (ARM)
$ws.Range("H32").Value = 9575.575999999999
$ws.Range("I32").Value = 7533.616
$ws.Range("J32").Value = 21997.5
$ws.Range("K32").Value = 7533.616
$ws.Range("L32").Value = 21997.5
$ws.Range("M32").Value = -7246.616
$ws.Range("N32").Value = -22571.5

# Row 45 (ARM)
$ws.Range("H45").Value = 931.5625
$ws.Range("I45").Value = 836.0714
$ws.Range("K45").Value = 836.0714
$ws.Range("M45").Value = -459.0714

# Row 116 (ARM)
$ws.Range("H116").Value = 1000
$ws.Range("I116").Value = 755.8823
$ws.Range("J116").Value = 1830
$ws.Range("K116").Value = 755.8823
$ws.Range("L116").Value = 1830
$ws.Range("M116").Value = 1538.1177
$ws.Range("N116").Value = -6418

# Row 122 (ARM)
$ws.Range("H122").Value = 2210.9167
$ws.Range("I122").Value = 2202.4
$ws.Range("J122").Value = 2253.5
$ws.Range("K122").Value = 6607.200000000001
$ws.Range("L122").Value = 6760.5
$ws.Range("M122").Value = -4157.200000000001
$ws.Range("N122").Value = -11660.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 755.8823
$ws.Range("J3").Value = 1830
$ws.Range("K3").Value = 755.8823
$ws.Range("L3").Value = 1830
$ws.Range("M3").Value = -641.8823
$ws.Range("N3").Value = -2058

# Row 99 (BSM)
$ws.Range("H99").Value = 40001044
$ws.Range("I99").Value = 50000950
$ws.Range("K99").Value = 50000950
$ws.Range("M99").Value = -49999452

# Row 105 (BSM)
$ws.Range("H105").Value = 50001228
$ws.Range("I105").Value = 76924160
$ws.Range("K105").Value = 76924160
$ws.Range("M105").Value = -76922413

# Row 107 (BSM)
$ws.Range("H107").Value = 868.56665
$ws.Range("I107").Value = 606.2083
$ws.Range("K107").Value = 606.2083
$ws.Range("M107").Value = 1313.7917

$ws = $wb.Worksheets.Item("CRP")
# Row 32 (CRP)
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# Row 43 (CRP)
$ws.Range("H43").Value = 11280.909
$ws.Range("J43").Value = 11280.909
$ws.Range("L43").Value = 11280.909
$ws.Range("N43").Value = -11648.909

# Row 55 (CRP)
$ws.Range("H55").Value = 8500
$ws.Range("I55").Value = 3000
$ws.Range("J55").Value = 14000
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 14000
$ws.Range("M55").Value = -2685
$ws.Range("N55").Value = -14630

# Row 101 (CRP)
$ws.Range("H101").Value = 11280.909
$ws.Range("J101").Value = 11280.909
$ws.Range("L101").Value = 11280.909
$ws.Range("N101").Value = -17770.909

$ws = $wb.Worksheets.Item("CUL")
# Row 23 (CUL)
$ws.Range("H23").Value = 350.33334
$ws.Range("J23").Value = 473.66666
$ws.Range("L23").Value = 1420.99998
$ws.Range("N23").Value = -1890.99998

# Row 70 (CUL)
$ws.Range("H70").Value = 11679.9375
$ws.Range("I70").Value = 35266.332
$ws.Range("J70").Value = 6236.923
$ws.Range("K70").Value = 105798.996
$ws.Range("L70").Value = 18710.769
$ws.Range("M70").Value = -105483.996
$ws.Range("N70").Value = -19340.769

# Row 73 (CUL)
$ws.Range("H73").Value = 11679.9375
$ws.Range("I73").Value = 35266.332
$ws.Range("J73").Value = 6236.923
$ws.Range("K73").Value = 105798.996
$ws.Range("L73").Value = 18710.769
$ws.Range("M73").Value = -104706.996
$ws.Range("N73").Value = -20894.769

# Row 87 (CUL)
$ws.Range("H87").Value = 1365.6666
$ws.Range("J87").Value = 1872.6666
$ws.Range("L87").Value = 5617.9998
$ws.Range("N87").Value = -8113.9998

# Row 88 (CUL)
$ws.Range("H88").Value = 5808.696
$ws.Range("J88").Value = 5808.696
$ws.Range("L88").Value = 17426.088
$ws.Range("N88").Value = -18282.088

# Row 90 (CUL)
$ws.Range("H90").Value = 1365.6666
$ws.Range("J90").Value = 1872.6666
$ws.Range("L90").Value = 16853.9994
$ws.Range("N90").Value = -29333.9994

# Row 91 (CUL)
$ws.Range("H91").Value = 5808.696
$ws.Range("J91").Value = 5808.696
$ws.Range("L91").Value = 17426.088
$ws.Range("N91").Value = -20390.088

# Row 109 (CUL)
$ws.Range("H109").Value = 7770958.5
$ws.Range("I109").Value = 10100936
$ws.Range("J109").Value = 4366.6665
$ws.Range("K109").Value = 30302808
$ws.Range("L109").Value = 13099.9995
$ws.Range("M109").Value = -30301768
$ws.Range("N109").Value = -15179.9995

# Row 131 (CUL)
$ws.Range("H131").Value = 25680230
$ws.Range("I131").Value = 83333760
$ws.Range("J131").Value = 56439.63
$ws.Range("K131").Value = 250001280
$ws.Range("L131").Value = 169318.89
$ws.Range("M131").Value = -249996240
$ws.Range("N131").Value = -179398.89

# Row 141 (CUL)
$ws.Range("H141").Value = 102502970
$ws.Range("I141").Value = 113890750
$ws.Range("K141").Value = 341672250
$ws.Range("M141").Value = -341667070

$ws = $wb.Worksheets.Item("GSM")
# Row 113 (GSM)
$ws.Range("H113").Value = 1220.1428
$ws.Range("I113").Value = 942.44446
$ws.Range("J113").Value = 1720
$ws.Range("K113").Value = 942.44446
$ws.Range("L113").Value = 1720
$ws.Range("M113").Value = 1227.55554
$ws.Range("N113").Value = -6060

# Row 132 (GSM)
$ws.Range("H132").Value = 3898.8572
$ws.Range("I132").Value = 3922.5
$ws.Range("K132").Value = 11767.5
$ws.Range("M132").Value = -9237.5

# Row 136 (GSM)
$ws.Range("H136").Value = 16735.934
$ws.Range("J136").Value = 16735.934
$ws.Range("L136").Value = 50207.802
$ws.Range("N136").Value = -55307.802

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW)
$ws.Range("H22").Value = 678.5714
$ws.Range("I22").Value = 902.75
$ws.Range("J22").Value = 588.9
$ws.Range("K22").Value = 902.75
$ws.Range("L22").Value = 588.9
$ws.Range("M22").Value = -607.75
$ws.Range("N22").Value = -1178.9

# Row 27 (LTW)
$ws.Range("H27").Value = 678.5714
$ws.Range("I27").Value = 902.75
$ws.Range("J27").Value = 588.9
$ws.Range("K27").Value = 902.75
$ws.Range("L27").Value = 588.9
$ws.Range("M27").Value = -795.75
$ws.Range("N27").Value = -802.9

# Row 132 (LTW)
$ws.Range("H132").Value = 2679.08
$ws.Range("I132").Value = 2425.2666
$ws.Range("K132").Value = 7275.7998
$ws.Range("M132").Value = -4745.7998

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (WVR)
$ws.Range("I62").Value = 100002200
$ws.Range("J62").Value = 3001.5
$ws.Range("K62").Value = 100002200
$ws.Range("L62").Value = 3001.5
$ws.Range("M62").Value = -100001576
$ws.Range("N62").Value = -4249.5

# Row 65 (WVR)
$ws.Range("I65").Value = 100002200
$ws.Range("J65").Value = 3001.5
$ws.Range("K65").Value = 500011000
$ws.Range("L65").Value = 15007.5
$ws.Range("M65").Value = -500007880
$ws.Range("N65").Value = -21247.5

# Row 126 (WVR)
$ws.Range("H126").Value = 83334270
$ws.Range("I126").Value = 90909930
$ws.Range("K126").Value = 272729790
$ws.Range("M126").Value = -272727320

# Row 136 (WVR)
$ws.Range("H136").Value = 1341.3914
$ws.Range("I136").Value = 942.8
$ws.Range("J136").Value = 2088.75
$ws.Range("K136").Value = 2828.4
$ws.Range("L136").Value = 6266.25
$ws.Range("M136").Value = -278.3999999999996
$ws.Range("N136").Value = -11366.25
